$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-17 Wednesday" "2024-07-18 Thursday"

Replace-Text "27×12=324" "39×57=2223"
Replace-Text "15×38=570" "75×73=5475"
Replace-Text "21×20=420" "45×15=675"
Replace-Text "37×27=999" "38×76=2888"
Replace-Text "38×59=2242" "20×42=840"

Replace-Text "62×64=3968" "96×58=5568"
Replace-Text "11×30=330" "52×80=4160"
Replace-Text "31×75=2325" "92×97=8924"
Replace-Text "65×89=5785" "73×15=1095"
Replace-Text "78×33=2574" "77×98=7546"

Replace-Text "32×32=1024" "48×46=2208"
Replace-Text "48×53=2544" "69×54=3726"
Replace-Text "61×39=2379" "52×20=1040"
Replace-Text "87×37=3219" "13×14=182"
Replace-Text "20×25=500" "29×28=812"

Replace-Text "36×20=720" "92×95=8740"
Replace-Text "41×38=1558" "47×87=4089"
Replace-Text "38×21=798" "89×58=5162"
Replace-Text "45×40=1800" "54×89=4806"
Replace-Text "66×40=2640" "86×68=5848"

Replace-Text "89×76=6764" "96×62=5952"
Replace-Text "14×57=798" "37×33=1221"
Replace-Text "87×21=1827" "15×27=405"
Replace-Text "39×95=3705" "17×42=714"
Replace-Text "63×49=3087" "50×19=950"
